# Apply updated odds values to Sheet1 as described in the commit diff.
# The workbook's single worksheet "Sheet1" contains match odds data; this
# script updates the specific cells that changed between the previous and
# current scrape of FlashScore odds for 2025-02-12 fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 6.5
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 1.48
$ws.Range("Y2").Value = 1.91
$ws.Range("Z2").Value = 1.91
$ws.Range("AA2").Value = 19
$ws.Range("AC2").Value = 21
$ws.Range("AK2").Value = 351
$ws.Range("AL2").Value = 8
$ws.Range("AM2").Value = 8

# Row 3
$ws.Range("H3").Value = 3
$ws.Range("K3").Value = 1.83
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25
$ws.Range("Q3").Value = 2.88
$ws.Range("R3").Value = 1.4
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 1.17
$ws.Range("U3").Value = 6
$ws.Range("V3").Value = 1.13
$ws.Range("Y3").Value = 2.25
$ws.Range("Z3").Value = 1.57
$ws.Range("AB3").Value = 10
$ws.Range("AJ3").Value = 81
$ws.Range("AL3").Value = 7
$ws.Range("AP3").Value = 34

# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 2.9
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 1.8
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5
$ws.Range("S5").Value = 5.6
$ws.Range("AA5").Value = 5
$ws.Range("AC5").Value = 11
$ws.Range("AG5").Value = 5.5
$ws.Range("AM5").Value = 17
$ws.Range("AN5").Value = 15
$ws.Range("AO5").Value = 41
$ws.Range("AR5").Value = 2.32
$ws.Range("AS5").Value = 1.62

# Row 6
$ws.Range("G6").Value = 1.95
$ws.Range("I6").Value = 4.5
$ws.Range("J6").Value = 2.75
$ws.Range("L6").Value = 5.5
$ws.Range("AB6").Value = 7.5
$ws.Range("AI6").Value = 23
$ws.Range("AM6").Value = 21
$ws.Range("AN6").Value = 17

# Row 7
$ws.Range("G7").Value = 7
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 1.55
$ws.Range("J7").Value = 7
$ws.Range("AN7").Value = 8.5

# Row 8
$ws.Range("G8").Value = 1.8
$ws.Range("I8").Value = 5.25
$ws.Range("J8").Value = 2.6
$ws.Range("AA8").Value = 5
$ws.Range("AB8").Value = 7
$ws.Range("AN8").Value = 19
$ws.Range("AO8").Value = 51

# Row 9
$ws.Range("G9").Value = 1.44
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 8.5
$ws.Range("J9").Value = 2.05
$ws.Range("L9").Value = 8.5
$ws.Range("Y9").Value = 2.5
$ws.Range("Z9").Value = 1.5
$ws.Range("AA9").Value = 5
$ws.Range("AD9").Value = 9
$ws.Range("AG9").Value = 7.5
$ws.Range("AH9").Value = 8
$ws.Range("AI9").Value = 26
$ws.Range("AJ9").Value = 101
$ws.Range("AL9").Value = 15
$ws.Range("AN9").Value = 26
$ws.Range("AO9").Value = 101
$ws.Range("AP9").Value = 67

# Row 10
$ws.Range("G10").Value = 2.8
$ws.Range("I10").Value = 2.63
$ws.Range("J10").Value = 3.75
$ws.Range("L10").Value = 3.5
$ws.Range("AB10").Value = 12
$ws.Range("AD10").Value = 29
$ws.Range("AO10").Value = 26
$ws.Range("AP10").Value = 26

# Row 13
$ws.Range("G13").Value = 2.5
$ws.Range("I13").Value = 3.25
$ws.Range("K13").Value = 1.88
$ws.Range("W13").Value = 1.53
$ws.Range("Y13").Value = 2.05
$ws.Range("AA13").Value = 6.2
$ws.Range("AB13").Value = 11
$ws.Range("AD13").Value = 28
$ws.Range("AL13").Value = 7.2
$ws.Range("AM13").Value = 15.5

# Row 14
$ws.Range("G14").Value = 1.33
$ws.Range("H14").Value = 4.45
$ws.Range("I14").Value = 10
$ws.Range("J14").Value = 1.85
$ws.Range("L14").Value = 8.5
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 7.3
$ws.Range("O14").Value = 1.3
$ws.Range("P14").Value = 3.2
$ws.Range("Q14").Value = 1.9
$ws.Range("R14").Value = 1.8
$ws.Range("U14").Value = 3.15
$ws.Range("V14").Value = 1.31
$ws.Range("X14").Value = 2.65
$ws.Range("Y14").Value = 2.37
$ws.Range("Z14").Value = 1.52
$ws.Range("AA14").Value = 5.4
$ws.Range("AB14").Value = 5.4
$ws.Range("AC14").Value = 9
$ws.Range("AD14").Value = 7.6
$ws.Range("AE14").Value = 12.5
$ws.Range("AF14").Value = 40
$ws.Range("AG14").Value = 7.3
$ws.Range("AH14").Value = 9
$ws.Range("AI14").Value = 27
$ws.Range("AL14").Value = 21
$ws.Range("AM14").Value = 70
$ws.Range("AN14").Value = 32
$ws.Range("AO14").Value = 350
$ws.Range("AP14").Value = 150

# Row 15
$ws.Range("G15").Value = 1.7
$ws.Range("I15").Value = 5
$ws.Range("U15").Value = 3.5
$ws.Range("V15").Value = 1.3
$ws.Range("AL15").Value = 13

# Row 16
$ws.Range("G16").Value = 1.7
$ws.Range("I16").Value = 5.5
$ws.Range("J16").Value = 2.4
$ws.Range("AA16").Value = 5.5
$ws.Range("AB16").Value = 7
$ws.Range("AD16").Value = 13
$ws.Range("AR16").Value = 1.83
$ws.Range("AS16").Value = 2.03

# Row 17
$ws.Range("G17").Value = 3.25
$ws.Range("I17").Value = 2.35
$ws.Range("J17").Value = 4
$ws.Range("L17").Value = 3.25
$ws.Range("O17").Value = 1.53
$ws.Range("P17").Value = 2.5
$ws.Range("AD17").Value = 34
$ws.Range("AG17").Value = 6.5
$ws.Range("AI17").Value = 19
$ws.Range("AM17").Value = 10
$ws.Range("AO17").Value = 23
$ws.Range("AR17").Value = 2.03
$ws.Range("AS17").Value = 1.83

# Row 18
$ws.Range("G18").Value = 1.55
$ws.Range("I18").Value = 6
$ws.Range("L18").Value = 6.5
$ws.Range("Y18").Value = 2.2
$ws.Range("Z18").Value = 1.62
$ws.Range("AA18").Value = 5.5
$ws.Range("AJ18").Value = 81

# Row 19
$ws.Range("G19").Value = 2.5
$ws.Range("H19").Value = 3.4
$ws.Range("I19").Value = 2.7
$ws.Range("L19").Value = 3.4
$ws.Range("U19").Value = 3.75
$ws.Range("V19").Value = 1.29
$ws.Range("AA19").Value = 8
$ws.Range("AD19").Value = 26
$ws.Range("AG19").Value = 9.5
$ws.Range("AK19").Value = 301
$ws.Range("AN19").Value = 10

# Row 20
$ws.Range("I20").Value = 3.2
$ws.Range("Q20").Value = 2.03
$ws.Range("R20").Value = 1.83
$ws.Range("U20").Value = 3.5
$ws.Range("V20").Value = 1.3

# Row 26
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 3.3
$ws.Range("I26").Value = 3.7
$ws.Range("J26").Value = 2.75
$ws.Range("L26").Value = 4.33
$ws.Range("O26").Value = 1.36
$ws.Range("P26").Value = 3
$ws.Range("AB26").Value = 9
$ws.Range("AD26").Value = 17
$ws.Range("AE26").Value = 17
$ws.Range("AH26").Value = 6.5
$ws.Range("AL26").Value = 10
$ws.Range("AP26").Value = 34
